# Refresh "cryptos" sheet data (GitHub Actions scheduled update).
# Updates Price (D) and Volume(1h) (E) columns, and for rows 44-48 also
# the Coin name (B) / Link (C) columns, which shifted due to re-ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.082.98"
$ws.Range("E2").Value = "  +0.16%  "
# Row 3
$ws.Range("D3").Value = "2.053.92"
$ws.Range("E3").Value = "  +0.03%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.35"
$ws.Range("E5").Value = "  +0.06%  "
# Row 6
$ws.Range("E6").Value = "  -2.09%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.06"
$ws.Range("E7").Value = "  +10.87%  "
# Row 8
$ws.Range("E8").Value = "  +0.00%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.386"
$ws.Range("E9").Value = "  +1.70%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0793"
$ws.Range("E10").Value = "  +0.53%  "
# Row 11
$ws.Range("E11").Value = "  +2.00%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.14"
$ws.Range("E12").Value = "  +8.77%  "
# Row 13
$ws.Range("D13").Value = "2.353.34"
$ws.Range("E13").Value = "  +0.06%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.817"
$ws.Range("E14").Value = "  +0.51%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.71"
$ws.Range("E15").Value = "  +9.22%  "
# Row 16
$ws.Range("D16").Value = "2.050.41"
$ws.Range("E16").Value = "  -0.14%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.78"
$ws.Range("E17").Value = "  +32.68%  "
# Row 18
$ws.Range("D18").Value = "37.061.17"
$ws.Range("E18").Value = "  +0.22%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "75.24"
$ws.Range("E19").Value = "  +3.70%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -2.28%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.42"
$ws.Range("E21").Value = "  +1.37%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.14"
$ws.Range("E22").Value = "  +0.97%  "
# Row 23
$ws.Range("E23").Value = "  +0.05%  "
# Row 24
$ws.Range("E24").Value = "  -0.05%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("E25").Value = "  +13.01%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.65"
$ws.Range("E26").Value = "  -0.85%  "
# Row 27
$ws.Range("E27").Value = "  +4.46%  "
# Row 28
$ws.Range("E28").Value = "  -0.05%  "
# Row 29
$ws.Range("E29").Value = "  -0.11%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  +10.72%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.86"
$ws.Range("E31").Value = "  +6.84%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0622"
$ws.Range("E32").Value = "  +0.45%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  +4.31%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0886"
$ws.Range("E34").Value = "  +2.92%  "
# Row 35
$ws.Range("E35").Value = "  -0.08%  "
# Row 36
$ws.Range("E36").Value = "  -0.85%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("E37").Value = "  -1.34%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +3.80%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.34"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.20"
$ws.Range("E40").Value = "  +28.23%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.11"
$ws.Range("E41").Value = "  +12.07%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.91"
$ws.Range("E42").Value = "  +0.48%  "
# Row 43
$ws.Range("E43").Value = "  +0.56%  "
# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.78"
$ws.Range("E44").Value = "  +1.70%  "
# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.14"
$ws.Range("E45").Value = "  -0.08%  "
# Row 46
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.06"
$ws.Range("E46").Value = "  -1.45%  "
# Row 47
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("E47").Value = "  +2.76%  "
# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.290.22"
$ws.Range("E48").Value = "  -0.10%  "
# Row 49
$ws.Range("E49").Value = "  -1.38%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.84"
$ws.Range("E50").Value = "  +0.61%  "
# Row 51
$ws.Range("D51").Value = "2.238.09"
$ws.Range("E51").Value = "  -0.05%  "
